$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header for new column E
$ws.Range("E1").Value = "vocabulary_concept_id"

# Values for E2:E59 -- sequential integers starting at 44819096
$startVal = 44819096
for ($row = 2; $row -le 59; $row++) {
    $ws.Cells.Item($row, 5).Value = $startVal + ($row - 2)
}

# Update view: scroll to A20, select A59
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 20
$ws.Range("A59").Select()
